$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 57, pushing the existing
# rows 57-61 down to 59-63.
$ws.Range("A57:A58").EntireRow.Insert()

# New row 57: Membrillo Champion "Especial" price update (week of 2023-04-18).
$ws.Range("A57").Value = 9
$ws.Range("B57").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C57").Value = "Metropolitana"
$ws.Range("D57").Value = 45034
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100104
$ws.Range("H57").Value = "Frutos de pepita"
$ws.Range("I57").Value = 100104003
$ws.Range("J57").Value = "Membrillo"
$ws.Range("K57").Value = "Champion"
$ws.Range("L57").Value = "Especial"
$ws.Range("M57").Value = 330
$ws.Range("N57").Value = 11000
$ws.Range("O57").Value = 11000
$ws.Range("P57").Value = 11000
$ws.Range("Q57").Value = '$/caja 18 kilos granel'
$ws.Range("R57").Value = "Región de O'Higgins"
$ws.Range("S57").Value = 611
$ws.Range("T57").Value = 18

# New row 58: Membrillo Champion "Primera" price update (week of 2023-04-18).
$ws.Range("A58").Value = 9
$ws.Range("B58").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C58").Value = "Metropolitana"
$ws.Range("D58").Value = 45034
$ws.Range("E58").Value = 13
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100104
$ws.Range("H58").Value = "Frutos de pepita"
$ws.Range("I58").Value = 100104003
$ws.Range("J58").Value = "Membrillo"
$ws.Range("K58").Value = "Champion"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 280
$ws.Range("N58").Value = 7500
$ws.Range("O58").Value = 7500
$ws.Range("P58").Value = 7500
$ws.Range("Q58").Value = '$/caja 18 kilos granel'
$ws.Range("R58").Value = "Región de O'Higgins"
$ws.Range("S58").Value = 417
$ws.Range("T58").Value = 18
